# feat: add 2022-Q1 data
#
# 1) Create a new "2022-Q1" worksheet (by duplicating the "2021-Q4" sheet so
#    the sheetPr/sheetFormatPr/pageMargins/header-row/column-A styles match
#    the other quarterly fund-holding sheets exactly), positioned right
#    before the "总计" (totals) sheet.
# 2) Fill the new sheet's data rows (2-9) with the 2022-Q1 fund holdings.
# 3) Update the "总计" sheet: insert the new 2022-Q1 summary row at the top
#    of the data and shift the existing quarters down by one, renumbering
#    the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: duplicate "2021-Q4" -> "2022-Q1", placed immediately after it
# (i.e. immediately before "总计").
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy([System.Reflection.Missing]::Value, $q4)
$q1new = $wb.Worksheets.Item("2021-Q4 (2)")
$q1new.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# Step 2: overwrite the data rows of the new sheet with the 2022-Q1 figures.
# Header row (B1:H1) is already identical, so it is left untouched.
# Columns B,D,E,F,G are text-like (fund code / ratios kept as text, exactly
# as authored in the source workbook) -- a leading apostrophe forces Excel
# to store them as text instead of auto-converting to numbers (this also
# protects fund codes with leading zeros, e.g. "001882"). Column C is plain
# text (fund name) and needs no such protection. Column H is a genuine
# number (ranking).
# ---------------------------------------------------------------------------

# Row 2
$q1new.Range("B2").Value = "'166005"
$q1new.Range("C2").Value = "中欧价值发现混合 -A"
$q1new.Range("D2").Value = "'43.52"
$q1new.Range("E2").Value = "'93.97"
$q1new.Range("F2").Value = "'4.55"
$q1new.Range("G2").Value = "'1.9802"
$q1new.Range("H2").Value = 3

# Row 3
$q1new.Range("B3").Value = "'001882"
$q1new.Range("C3").Value = "中欧价值发现混合 -E"
$q1new.Range("D3").Value = "'43.52"
$q1new.Range("E3").Value = "'93.97"
$q1new.Range("F3").Value = "'4.55"
$q1new.Range("G3").Value = "'1.9802"
$q1new.Range("H3").Value = 3

# Row 4
$q1new.Range("B4").Value = "'001810"
$q1new.Range("C4").Value = "中欧潜力价值灵活配置混合A"
$q1new.Range("D4").Value = "'28.67"
$q1new.Range("E4").Value = "'94.05"
$q1new.Range("F4").Value = "'4.59"
$q1new.Range("G4").Value = "'1.3160"
$q1new.Range("H4").Value = 3

# Row 5
$q1new.Range("B5").Value = "'004232"
$q1new.Range("C5").Value = "中欧价值发现混合 -C"
$q1new.Range("D5").Value = "'10.98"
$q1new.Range("E5").Value = "'93.97"
$q1new.Range("F5").Value = "'4.55"
$q1new.Range("G5").Value = "'0.4996"
$q1new.Range("H5").Value = 3

# Row 6
$q1new.Range("B6").Value = "'166024"
$q1new.Range("C6").Value = "中欧恒利三年定期开放混合"
$q1new.Range("D6").Value = "'4.48"
$q1new.Range("E6").Value = "'98.71"
$q1new.Range("F6").Value = "'4.94"
$q1new.Range("G6").Value = "'0.2213"
$q1new.Range("H6").Value = 4

# Row 7
$q1new.Range("B7").Value = "'005764"
$q1new.Range("C7").Value = "中欧潜力价值灵活配置混合C"
$q1new.Range("D7").Value = "'3.43"
$q1new.Range("E7").Value = "'94.05"
$q1new.Range("F7").Value = "'4.59"
$q1new.Range("G7").Value = "'0.1574"
$q1new.Range("H7").Value = 3

# Row 8
$q1new.Range("B8").Value = "'001891"
$q1new.Range("C8").Value = "中欧成长优选回报灵活配置混合E"
$q1new.Range("D8").Value = "'2.97"
$q1new.Range("E8").Value = "'94.42"
$q1new.Range("F8").Value = "'3.75"
$q1new.Range("G8").Value = "'0.1114"
$q1new.Range("H8").Value = 4

# Row 9
$q1new.Range("B9").Value = "'166020"
$q1new.Range("C9").Value = "中欧成长优选回报灵活配置混合A"
$q1new.Range("D9").Value = "'2.97"
$q1new.Range("E9").Value = "'94.42"
$q1new.Range("F9").Value = "'3.75"
$q1new.Range("G9").Value = "'0.1114"
$q1new.Range("H9").Value = 4

# ---------------------------------------------------------------------------
# Step 3: update the "总计" (totals) sheet. Shift the five existing quarters
# down one row (re-numbering the leading index column A), then write the new
# 2022-Q1 totals into row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# New row created at the bottom needs the same column-A style ("s=2") as the
# rest of the index column, so copy formatting down from the row above it
# before writing its value.
$total.Range("A6").Copy($total.Range("A7"))

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 6
$total.Range("D7").Value = 0.25

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.07000000000000001

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 5
$total.Range("D5").Value = 0.8100000000000001

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 8
$total.Range("D4").Value = 2.77

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 8
$total.Range("D3").Value = 3.1

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 6.38

# ---------------------------------------------------------------------------
# Restore the original active sheet/selection (the sheet-copy operation
# above left the new "2022-Q1" sheet selected; the source workbook had the
# first sheet, "2020-Q4", active).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
